$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 10's formatting down into the new row 11 so the new
# data row picks up the same cell styles (font/fill/border) used
# by the rest of the imported data rows.
$ws.Range("A10:S10").Copy()
$ws.Range("A11:S11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows("11").RowHeight = 15

# Populate the new placement row (Yafang Deng / GMC 7000019) that
# reproduces the StringConverter test-data issue.
$ws.Range("A11").Value = "Yafang"
$ws.Range("B11").Value = "Deng"
$ws.Range("C11").Value = 7000019
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = "WMD/5AT01/085/ST1/001"
$ws.Range("G11").Value = "01/01/2019"
$ws.Range("H11").Value = "01/01/2026"
$ws.Range("I11").Value = "In Post"
$ws.Range("J11").Value = "CURRENT"
$ws.Range("K11").Value = "This is for ""test\"
$ws.Range("L11").Value = "1"
$ws.Range("M11").Value = "Specialty Training Year 1"
$ws.Range("N11").Value = ""
$ws.Range("O11").Value = ""
$ws.Range("P11").Value = ""
$ws.Range("Q11").Value = ""
$ws.Range("R11").Value = ""
$ws.Range("S11").Value = ""

[void]$ws.Range("A11").Select()
